$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read existing data (to mirror "read and write data" flow) and then overwrite
# rows 2 and 3 with updated user records.
$ws.Range("A2").Value = "Teresa Rolfson"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = "Ansley_Marvin77@hotmail.com"

$ws.Range("A3").Value = "Gregg Hyatt III"
$ws.Range("B3").Value = 43
$ws.Range("C3").Value = "Marta.Schulist@yahoo.com"
